$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022115658001729
$ws.Range("D2").Value = 1.027273437115517
$ws.Range("E2").Value = 1.022911517641687
$ws.Range("I2").Value = 1.031662037948412
$ws.Range("J2").Value = 1.027303309300051
$ws.Range("K2").Value = 1.030093645979396
$ws.Range("L2").Value = 1.025744504756261
$ws.Range("N2").Value = 1.013192062822364

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022820545736574
$ws.Range("D3").Value = 1.027777706988114
$ws.Range("E3").Value = 1.023503056282056
$ws.Range("I3").Value = 1.031787925071001
$ws.Range("J3").Value = 1.027647379312717
$ws.Range("K3").Value = 1.030406427222679
$ws.Range("L3").Value = 1.026143397959211
$ws.Range("N3").Value = 1.013305175440048

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0232771354193
$ws.Range("D4").Value = 1.028104281836498
$ws.Range("E4").Value = 1.023886641224532
$ws.Range("I4").Value = 1.031868170550839
$ws.Range("J4").Value = 1.027869797007596
$ws.Range("K4").Value = 1.03060837936504
$ws.Range("L4").Value = 1.026401630347836
$ws.Range("N4").Value = 1.013378292950636

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023469198572965
$ws.Range("D5").Value = 1.028241638560583
$ws.Range("E5").Value = 1.024048094521317
$ws.Range("I5").Value = 1.031901614693397
$ws.Range("J5").Value = 1.02796324768325
$ws.Range("K5").Value = 1.030693173720757
$ws.Range("L5").Value = 1.026510218874879
$ws.Range("N5").Value = 1.013409013390558

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02350145337795
$ws.Range("D6").Value = 1.028264705091632
$ws.Range("E6").Value = 1.024075214566434
$ws.Range("I6").Value = 1.031907213017467
$ws.Range("J6").Value = 1.027978935272132
$ws.Range("K6").Value = 1.03070740479018
$ws.Range("L6").Value = 1.026528452942059
$ws.Range("N6").Value = 1.013414170409031

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02327970133578
$ws.Range("D7").Value = 1.02810611695226
$ws.Range("E7").Value = 1.023888797811011
$ws.Range("I7").Value = 1.031868618577924
$ws.Range("J7").Value = 1.027871045912942
$ws.Range("K7").Value = 1.030609512810932
$ws.Range("L7").Value = 1.026403081206411
$ws.Range("N7").Value = 1.013378703510607

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.022353777421229
$ws.Range("D8").Value = 1.027443798597268
$ws.Range("E8").Value = 1.023111259873501
$ws.Range("I8").Value = 1.031704832468662
$ws.Range("J8").Value = 1.027419633803284
$ws.Range("K8").Value = 1.030199441758152
$ws.Range("L8").Value = 1.025879286440023
$ws.Range("N8").Value = 1.013230304786373

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020725951910421
$ws.Range("D9").Value = 1.026278931342891
$ws.Range("E9").Value = 1.021747510562603
$ws.Range("I9").Value = 1.031406979503315
$ws.Range("J9").Value = 1.026622574956744
$ws.Range("K9").Value = 1.029473550043155
$ws.Range("L9").Value = 1.024957293693395
$ws.Range("N9").Value = 1.012968261719713

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.01964339594707
$ws.Range("D10").Value = 1.025503967944321
$ws.Range("E10").Value = 1.020842753661685
$ws.Range("I10").Value = 1.031202252968525
$ws.Range("J10").Value = 1.026090200692033
$ws.Range("K10").Value = 1.02898749159129
$ws.Range("L10").Value = 1.024343396023542
$ws.Range("N10").Value = 1.012793227331794

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.019175295071619
$ws.Range("D11").Value = 1.025168808823156
$ws.Range("E11").Value = 1.020452056347719
$ws.Range("I11").Value = 1.031112155211243
$ws.Range("J11").Value = 1.025859456707868
$ws.Range("K11").Value = 1.028776535375593
$ws.Range("L11").Value = 1.024077772192803
$ws.Range("N11").Value = 1.012717360986333

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019001521833355
$ws.Range("D12").Value = 1.025044378756467
$ws.Range("E12").Value = 1.020307096658407
$ws.Range("I12").Value = 1.03107847200322
$ws.Range("J12").Value = 1.025773716123834
$ws.Range("K12").Value = 1.028698104713994
$ws.Range("L12").Value = 1.023979138986662
$ws.Range("N12").Value = 1.01268917001149

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019038792210758
$ws.Range("D13").Value = 1.025071066545958
$ws.Range("E13").Value = 1.020338183606984
$ws.Range("I13").Value = 1.031085706965475
$ws.Range("J13").Value = 1.025792109195204
$ws.Range("K13").Value = 1.028714931598091
$ws.Range("L13").Value = 1.024000294702005
$ws.Range("N13").Value = 1.012695217554008

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.019160928863577
$ws.Range("D14").Value = 1.025158522098755
$ws.Range("E14").Value = 1.020440070603227
$ws.Range("I14").Value = 1.031109375365868
$ws.Range("J14").Value = 1.025852370008918
$ws.Range("K14").Value = 1.028770053735753
$ws.Range("L14").Value = 1.024069618493364
$ws.Range("N14").Value = 1.012715030929861

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.01923619463768
$ws.Range("D15").Value = 1.025212414755476
$ws.Range("E15").Value = 1.020502868167085
$ws.Range("I15").Value = 1.031123929536335
$ws.Range("J15").Value = 1.025889494487967
$ws.Range("K15").Value = 1.028804006789111
$ws.Range("L15").Value = 1.024112335365375
$ws.Range("N15").Value = 1.012727237182986

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.01967447620319
$ws.Range("D16").Value = 1.025526220065492
$ws.Range("E16").Value = 1.020868705651323
$ws.Range("I16").Value = 1.031208201993315
$ws.Range("J16").Value = 1.026105509845672
$ws.Range("K16").Value = 1.029001481866303
$ws.Range("L16").Value = 1.024361028906908
$ws.Range("N16").Value = 1.012798260785523

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019949575004522
$ws.Range("D17").Value = 1.025723171789079
$ws.Range("E17").Value = 1.021098473416854
$ws.Range("I17").Value = 1.031260676378306
$ws.Range("J17").Value = 1.026240952029963
$ws.Range("K17").Value = 1.029125222677576
$ws.Range("L17").Value = 1.024517081896698
$ws.Range("N17").Value = 1.012842792191945

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.020110098341184
$ws.Range("D18").Value = 1.025838089348399
$ws.Range("E18").Value = 1.02123259600626
$ws.Range("I18").Value = 1.031291143850731
$ws.Range("J18").Value = 1.02631993160971
$ws.Range("K18").Value = 1.02919735125401
$ws.Range("L18").Value = 1.024608123955476
$ws.Range("N18").Value = 1.012868759310594

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.020164843247016
$ws.Range("D19").Value = 1.025877279809617
$ws.Range("E19").Value = 1.021278345737207
$ws.Range("I19").Value = 1.031301508705491
$ws.Range("J19").Value = 1.026346857906032
$ws.Range("K19").Value = 1.029221937160117
$ws.Range("L19").Value = 1.024639170117759
$ws.Range("N19").Value = 1.012877612174726

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01992005298913
$ws.Range("D20").Value = 1.025702036689129
$ws.Range("E20").Value = 1.021073810846951
$ws.Range("I20").Value = 1.031255060842243
$ws.Range("J20").Value = 1.026226422584672
$ws.Range("K20").Value = 1.029111951346342
$ws.Range("L20").Value = 1.024500336912157
$ws.Range("N20").Value = 1.012838015144832

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.019124959878221
$ws.Range("D21").Value = 1.025132766863683
$ws.Range("E21").Value = 1.020410062913861
$ws.Range("I21").Value = 1.031102411591622
$ws.Range("J21").Value = 1.025834625573052
$ws.Range("K21").Value = 1.028753823619688
$ws.Range("L21").Value = 1.024049203487662
$ws.Range("N21").Value = 1.012709196679816

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018625634477147
$ws.Range("D22").Value = 1.02477520958866
$ws.Range("E22").Value = 1.019993680583264
$ws.Range("I22").Value = 1.031005180553325
$ws.Range("J22").Value = 1.025588102832205
$ws.Range("K22").Value = 1.02852823789101
$ws.Range("L22").Value = 1.023765739934414
$ws.Range("N22").Value = 1.012628140938286

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018890280507571
$ws.Range("D23").Value = 1.024964722148683
$ws.Range("E23").Value = 1.020214322687797
$ws.Range("I23").Value = 1.031056843158629
$ws.Range("J23").Value = 1.025718806207868
$ws.Range("K23").Value = 1.028647864163836
$ws.Range("L23").Value = 1.023915991556305
$ws.Range("N23").Value = 1.012671115882085

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019933392523188
$ws.Range("D24").Value = 1.025711586610979
$ws.Range("E24").Value = 1.021084954482937
$ws.Range("I24").Value = 1.031257598693812
$ws.Range("J24").Value = 1.026232987882378
$ws.Range("K24").Value = 1.029117948236151
$ws.Range("L24").Value = 1.02450790319143
$ws.Range("N24").Value = 1.012840173709409

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021146324000191
$ws.Range("D25").Value = 1.026579801640254
$ws.Range("E25").Value = 1.022099304421396
$ws.Range("I25").Value = 1.031485071041037
$ws.Range("J25").Value = 1.026828816381276
$ws.Range("K25").Value = 1.029661592137894
$ws.Range("L25").Value = 1.025195522909045
$ws.Range("N25").Value = 1.013036067916457
